$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full 3x3 sending/target cluster grid (ECs added) with refreshed Ptprk-Ptprk stats
$data = @(
    @("ECs", "Ptprk", "Ptprk", "ECs", 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 373.1673516515999, 3358.5061648644, 0.1906758458538097, 0.1906758458538098),
    @("ECs", "Ptprk", "Ptprk", "FAPs", 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 89.69618692253999, 807.2656823028599, 0.04583170589715605, 0.04583170589715606),
    @("ECs", "Ptprk", "Ptprk", "sCs", 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 391.7225120571399, 3525.502608514259, 0.2001569027834166, 0.2001569027834167),
    @("FAPs", "Ptprk", "Ptprk", "ECs", 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 89.69618692253999, 807.2656823028599, 0.04583170589715605, 0.04583170589715606),
    @("FAPs", "Ptprk", "Ptprk", "FAPs", 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 21.559779849001, 194.038018641009, 0.01101631544382337, 0.01101631544382337),
    @("FAPs", "Ptprk", "Ptprk", "sCs", 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 94.15618892632433, 847.4057003369189, 0.04811061548240626, 0.04811061548240626),
    @("sCs", "Ptprk", "Ptprk", "ECs", 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 2, 0.6666666666666666, 19.31754, 57.95262, 0.4366644545343824, 0.4366644545343825, 391.7225120571399, 3525.502608514259, 0.2001569027834166, 0.2001569027834167),
    @("sCs", "Ptprk", "Ptprk", "FAPs", 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 3, 1, 4.643251, 13.929753, 0.1049586368233857, 0.1049586368233857, 94.15618892632433, 847.4057003369189, 0.04811061548240626, 0.04811061548240626),
    @("sCs", "Ptprk", "Ptprk", "sCs", 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 3, 1, 20.27807433333333, 60.83422299999999, 0.4583769086422319, 0.4583769086422319, 411.2002986681921, 3700.802688013728, 0.210109390376409, 0.210109390376409),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

Write-Output "done"
